# Apply the "Dry Cakes" (sheet2) additions + active-tab/selection change
# described in the diff.

$wb = $excel.ActiveWorkbook
$dryCakes = $wb.Worksheets.Item("Dry Cakes")

# --- New cheesecake rows (10-14) on the "Dry Cakes" sheet -----------------

# Column B (name) - set first, in row order.
$dryCakes.Range("B10").Value = "New York Cheesecake"
$dryCakes.Range("B11").Value = "Lemon  Cheesecake"
$dryCakes.Range("B12").Value = "Strawberry  Cheesecake"
$dryCakes.Range("B13").Value = "Blueberry  Cheesecake"
$dryCakes.Range("B14").Value = "Biscoff  Cheesecake"

# Column A (id) - apply the same cell formatting used on the existing data
# rows (e.g. row 2) before filling in the values, in row order.
$dryCakes.Range("A2").Copy()
$dryCakes.Range("A10:A14").PasteSpecial(-4122)
$dryCakes.Range("A10").Value = "dc9"
$dryCakes.Range("A11").Value = "dc10"
$dryCakes.Range("A12").Value = "dc11"
$dryCakes.Range("A13").Value = "dc12"
$dryCakes.Range("A14").Value = "dc13"

# Column C (image) - note the source order is not simple row order.
$dryCakes.Range("C11").Value = "dry-cakes/Lemon-Cheesecake.jpg"
$dryCakes.Range("C12").Value = "dry-cakes/Strawberry-Cheesecake.jpg"
$dryCakes.Range("C13").Value = "dry-cakes/Blueberry-Cheesecake.jpg"
$dryCakes.Range("C14").Value = "dry-cakes/Biscoff-Cheesecake.jpg"
$dryCakes.Range("C10").Value = "dry-cakes/New-York-Cheesecake.jpg"

# Column F (description).
$dryCakes.Range("F10").Value = "[Veg preparation]"
$dryCakes.Range("F11").Value = "[Veg preparation]"
$dryCakes.Range("F12").Value = "[Veg preparation]"
$dryCakes.Range("F13").Value = "[Veg preparation]"
$dryCakes.Range("F14").Value = "[Veg preparation]"

# Column G (inStock) - apply formatting from row 2, then set the value.
$dryCakes.Range("G2").Copy()
$dryCakes.Range("G10:G14").PasteSpecial(-4122)
$dryCakes.Range("G10:G14").Value = "yes"

# Column H (onDiscount) - apply formatting from row 2, then set the value.
$dryCakes.Range("H2").Copy()
$dryCakes.Range("H10:H14").PasteSpecial(-4122)
$dryCakes.Range("H10:H14").Value = "no"

# --- Active tab / selection ------------------------------------------------
# Move the active tab / selection from "Sweet Delights" to "Dry Cakes".
$dryCakes.Activate()
[void]$dryCakes.Range("H15").Select()
